$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38
$ws.Range("B38").Value = 6781354
$ws.Range("F38").Value = "Puntarenas"
$ws.Range("G38").Value = "AD San Carlos"
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = "H"
$ws.Range("K38").Value = 2.4
$ws.Range("L38").Value = 3.2
$ws.Range("M38").Value = 2.8
$ws.Range("N38").Value = 2.3
$ws.Range("O38").Value = 3.2
$ws.Range("P38").Value = 3
$ws.Range("R38").Value = 2
$ws.Range("S38").Value = 1.8
$ws.Range("T38").Value = 2.25
$ws.Range("W38").Value = 1.3
$ws.Range("Y38").Value = -1
$ws.Range("Z38").Value = 1
$ws.Range("AA38").Value = -1
$ws.Range("AB38").Value = -1
$ws.Range("AC38").Value = 0.8999999999999999

# Row 39
$ws.Range("B39").Value = 6782522
$ws.Range("F39").Value = "Municipal Perez Zeledon"
$ws.Range("G39").Value = "Sporting San Jose"
$ws.Range("I39").Value = 2
$ws.Range("J39").Value = "A"
$ws.Range("K39").Value = 2.5
$ws.Range("L39").Value = 3.5
$ws.Range("M39").Value = 2.5
$ws.Range("N39").Value = 2.2
$ws.Range("O39").Value = 3.5
$ws.Range("P39").Value = 2.9
$ws.Range("R39").Value = 1.9
$ws.Range("S39").Value = 1.9
$ws.Range("T39").Value = 2.5
$ws.Range("W39").Value = -1
$ws.Range("Y39").Value = 1.9
$ws.Range("Z39").Value = -1
$ws.Range("AA39").Value = 0.8999999999999999
$ws.Range("AB39").Value = 0.8999999999999999
$ws.Range("AC39").Value = -1

# Row 224
$ws.Range("B224").Value = 7623946
$ws.Range("F224").Value = "Cartagines"
$ws.Range("G224").Value = "Alajuelense"
$ws.Range("K224").Value = 3.4
$ws.Range("L224").Value = 3.4
$ws.Range("M224").Value = 1.95
$ws.Range("N224").Value = 3.8
$ws.Range("O224").Value = 3.6
$ws.Range("P224").Value = 1.8
$ws.Range("Q224").Value = 0.5
$ws.Range("R224").Value = 2
$ws.Range("S224").Value = 1.8
$ws.Range("T224").Value = 2.75
$ws.Range("U224").Value = 1.975
$ws.Range("V224").Value = 1.825

# Row 225
$ws.Range("B225").Value = 7623944
$ws.Range("F225").Value = "Santos de Gupiles"
$ws.Range("G225").Value = "Municipal Liberia"
$ws.Range("K225").Value = 2.9
$ws.Range("L225").Value = 3.25
$ws.Range("M225").Value = 2.375
$ws.Range("N225").Value = 3.5
$ws.Range("O225").Value = 3.3
$ws.Range("P225").Value = 2.1
$ws.Range("R225").Value = 1.95
$ws.Range("S225").Value = 1.85
$ws.Range("U225").Value = 1.95
$ws.Range("V225").Value = 1.85

# Row 226
$ws.Range("R226").Value = 1.85
$ws.Range("S226").Value = 1.95
$ws.Range("U226").Value = 2
$ws.Range("V226").Value = 1.8

# Row 227
$ws.Range("U227").Value = 1.9
$ws.Range("V227").Value = 1.9

# Row 229
$ws.Range("N229").Value = 1.285
$ws.Range("O229").Value = 5.25
$ws.Range("P229").Value = 8
$ws.Range("R229").Value = 1.875
$ws.Range("S229").Value = 1.925
$ws.Range("U229").Value = 1.825
$ws.Range("V229").Value = 1.975
